# Assessments.xlsx — "Added a couple more courseworks"
#
# 1. Row 19 (PH167 Online Assessment 1) loses its G (duration) value.
# 2. Row 21 (EE270 Exam) gains a G (duration) value of 1.5.
# 3. New rows 22-31 are appended: EE270 lab/logbook/assignment items and
#    a new EE312 module (quizzes, class test, exam).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: remove the duration (G19) value ---
$ws.Range("G19").Clear()

# --- Row 21: add a duration (G21) value ---
$ws.Range("G21").Value = 1.5

# --- New data rows for EE270 (labs / logbook / assignment) ---
$newRows = @(
    @{ Row = 22; A = "EE270"; B = "C"; C = 1; D = 5;   E = 10; F = "Lab 1 " },
    @{ Row = 23; A = "EE270"; B = "C"; C = 1; D = 6;   E = 10; F = "Lab 2" },
    @{ Row = 24; A = "EE270"; B = "C"; C = 1; D = 7;   E = 5;  F = "Logbook" },
    @{ Row = 25; A = "EE270"; B = "C"; C = 1; D = 7;   E = 15; F = "Assignment" },
    @{ Row = 26; A = "EE312"; B = "C"; C = 1; D = 7;   E = 5;  F = "Quiz 1" },
    @{ Row = 27; A = "EE312"; B = "C"; C = 1; D = 9;   E = 5;  F = "Quiz 2" },
    @{ Row = 28; A = "EE312"; B = "C"; C = 1; D = 11;  E = 5;  F = "Quiz 3" },
    @{ Row = 29; A = "EE312"; B = "C"; C = 1; D = "E"; E = 15; F = "Microcontroller Demo" },
    @{ Row = 30; A = "EE312"; B = "C"; C = 2; D = 7;   E = 10; F = "Class Test" }
)

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
}

# --- New row 31: EE312 Exam, with a duration value ---
$ws.Range("A31").Value = "EE312"
$ws.Range("B31").Value = "E"
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = "E"
$ws.Range("E31").Value = 60
$ws.Range("F31").Value = "Exam"
$ws.Range("G31").Value = 2

# --- Move the active selection to match the saved view state ---
$ws.Range("F34").Select()
